$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Change zoom value from 20 to 14
$ws.Range("B4").Value = 14

# Change label at A12 from "min zoom:" to "benchmark zoom:"
$ws.Range("A12").Value = "benchmark zoom:"

# Update B13 formula (new scaling formula, no longer clamping to min zoom)
$ws.Range("B13").Formula = "=`$B11/POWER(2,B4-B12)"

# Add empty formatted cell at C13 (matches the formatting used by the
# other input-row placeholder cells such as C3/E3)
$ws.Range("C13").HorizontalAlignment = -4152

# Update selection to I14
$ws.Range("I14").Select()
